# Add a new worksheet "ciudad_distinta" at the end of the workbook,
# populate it with the distinct city values, and make it the active sheet.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "ciudad_distinta"

$newSheet.Range("A1").Value = "ciudad"
$newSheet.Range("A2").Value = "Ciudad8"
$newSheet.Range("A3").Value = "Ciudad1"

# Match the authored selection/active-cell state for the new sheet.
$newSheet.Range("B8").Select()
